# "Factor demand-distance fixed. Numbers of iterations as variable."
#
# Add a new column P to the Results_Summary sheet that computes, for each
# data row, the "demand/distance" factor: 50000 divided by (2 * Customers),
# where Customers lives in column D. This extends the used range from
# A1:O28 to A1:P28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 3
$lastDataRow = 28
$targetColumn = 16   # column P

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $cell = $ws.Cells.Item($row, $targetColumn)
    $cell.Formula = "=50000/(D$row*2)"
}

# Nudge the default column width the same way the source workbook records it.
$ws.StandardWidth = 11.55078125

# Leave the active selection on C12, matching the saved workbook state.
$ws.Range("C12").Select()
